$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new entry mirroring row 3's person/work item/status, with a new
# "next week" plan and a new completion-target description.
$ws.Range("A6").Value = [DateTime]"2019-12-02"
$ws.Range("B6").Value = "艾永芮"
$ws.Range("C6").Value = "資料蒐集與資料庫建置"
$ws.Range("D6").Value = "建置資料庫存放價量資料及模擬結果"
$ws.Range("E6").Value = "研究價量模擬方法"
$ws.Range("F6").Value = "支援價量模擬工作，並完成相應的資料庫建置"

$ws.Range("D6").WrapText = $true
$ws.Range("F6").WrapText = $true

$ws.Rows.Item(6).RowHeight = 57
$ws.Rows.Item(2).RowHeight = 76
$ws.Rows.Item(3).RowHeight = 95

$ws.Range("A6").Select() | Out-Null
